# Updated symbol list on Tue Dec 13 19:27:52 UTC 2022 with GitHub Actions
#
# Applies the per-cell value updates to the "cryptos" sheet. The D column
# holds numeric-looking values that are stored as TEXT (inline strings) in
# the workbook, so a plain Range.Value assignment would get auto-coerced to
# a real number by Excel. To preserve the original text type (and keep the
# cell's style untouched) each D-column write temporarily formats the cell
# as Text ("@"), assigns the literal string, then restores the cell style
# back to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

# --- Column D (Price) numeric-looking text updates ---
Set-TextValue "D2"  "268.98"
Set-TextValue "D3"  "22.79"
Set-TextValue "D4"  "6.333"
Set-TextValue "D6"  "3.646"
Set-TextValue "D7"  "6.663"
Set-TextValue "D8"  "1.388"
Set-TextValue "D9"  "0.8295"
Set-TextValue "D10" "0.01371"
Set-TextValue "D11" "0.1608"
Set-TextValue "D12" "0.08300"
Set-TextValue "D13" "0.03546"
Set-TextValue "D14" "0.03201"
Set-TextValue "D15" "0.09326"
Set-TextValue "D16" "3.838"
Set-TextValue "D17" "0.001647"
Set-TextValue "D18" "0.04744"
Set-TextValue "D19" "0.006345"
Set-TextValue "D20" "0.005658"
Set-TextValue "D21" "0.001077"
Set-TextValue "D23" "3.723"
Set-TextValue "D26" "0.1239"
Set-TextValue "D27" "0.0002705"
Set-TextValue "D40" "0.04718"
Set-TextValue "D41" "0.006983"

# --- Rows 42/43: CEJI and BKEXToken swap places (with new prices) ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1160"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003300"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

# --- remaining column D updates ---
Set-TextValue "D44" "0.01181"
Set-TextValue "D45" "0.00006254"
Set-TextValue "D46" "0.0009903"
Set-TextValue "D48" "0.9204"
Set-TextValue "D49" "0.002364"
Set-TextValue "D50" "0.00001401"
